# Updated cryptos list with GitHub Actions
# Applies per-cell value updates to the crypto price table on Sheet1,
# matching the upstream data refresh described in the commit diff.
#
# Cells in column D that look like plain numbers (e.g. "4.11") are
# written with a leading apostrophe so Excel stores them as text
# (quote-prefixed), exactly as they were originally authored (the
# sheet stores every Price/Volume value as text, including ones that
# are numeric-looking, so dotted "thousands" values like "27.009.67"
# keep rendering correctly instead of being parsed as numbers/dates).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)


# Row 2
$ws.Cells.Item(2, 4).Value = "27.009.67"
$ws.Cells.Item(2, 5).Value = "  +0.40%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.682.58"
$ws.Cells.Item(3, 5).Value = "  +0.54%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'215.68"
$ws.Cells.Item(5, 5).Value = "  -0.14%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.516"
$ws.Cells.Item(6, 5).Value = "  -2.63%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.12%  "

# Row 8
$ws.Cells.Item(8, 2).Value = "Cardano"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(8, 4).Value = "'0.251"
$ws.Cells.Item(8, 5).Value = "  -1.60%  "

# Row 9
$ws.Cells.Item(9, 2).Value = "Solana"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(9, 4).Value = "'21.27"
$ws.Cells.Item(9, 5).Value = "  +4.54%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -0.13%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.32%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "1.920.25"
$ws.Cells.Item(12, 5).Value = "  +0.53%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.720.59"
$ws.Cells.Item(13, 5).Value = "  +2.63%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'4.11"
$ws.Cells.Item(14, 5).Value = "  +0.29%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.533"
$ws.Cells.Item(15, 5).Value = "  +1.92%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'66.01"
$ws.Cells.Item(16, 5).Value = "  +0.37%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "27.024.58"
$ws.Cells.Item(17, 5).Value = "  +0.27%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'8.16"
$ws.Cells.Item(18, 5).Value = "  +3.85%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'236.16"
$ws.Cells.Item(19, 5).Value = "  +1.25%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.17%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +0.03%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -0.50%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -4.03%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'146.61"
$ws.Cells.Item(25, 5).Value = "  +0.53%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'7.23"
$ws.Cells.Item(26, 5).Value = "  +1.07%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'16.06"
$ws.Cells.Item(27, 5).Value = "  +0.53%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'0.112"
$ws.Cells.Item(28, 5).Value = "  -2.97%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.18%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'0.0500"
$ws.Cells.Item(30, 5).Value = "  +0.44%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -0.45%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.23%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "1.511.02"
$ws.Cells.Item(33, 5).Value = "  +3.25%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +0.32%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'1.68"
$ws.Cells.Item(35, 5).Value = "  +3.94%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -0.47%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +2.87%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.916"
$ws.Cells.Item(38, 5).Value = "  +0.99%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +3.18%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +6.79%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'5.74"
$ws.Cells.Item(41, 5).Value = "  -5.08%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -0.01%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'68.08"
$ws.Cells.Item(43, 5).Value = "  +3.30%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -1.12%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "1.826.11"
$ws.Cells.Item(45, 5).Value = "  +0.25%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.782"
$ws.Cells.Item(46, 5).Value = "  +0.16%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'90.23"
$ws.Cells.Item(47, 5).Value = "  -0.52%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "Algorand"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(48, 4).Value = "'0.104"
$ws.Cells.Item(48, 5).Value = "  +3.92%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "RenderToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(49, 4).Value = "'1.52"
$ws.Cells.Item(49, 5).Value = "  -0.99%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'7.81"
$ws.Cells.Item(50, 5).Value = "  +2.40%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  +0.08%  "

